$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 22754764
$ws.Range("J40").Value = 55607424
$ws.Range("L40").Value = 55607424
$ws.Range("N40").Value = -55607774
$ws.Range("H43").Value = 3080872.2
$ws.Range("J43").Value = 4915.3335
$ws.Range("L43").Value = 4915.3335
$ws.Range("N43").Value = -5053.3335
$ws.Range("H64").Value = 28860348
$ws.Range("J64").Value = 62504380
$ws.Range("L64").Value = 62504380
$ws.Range("N64").Value = -62504876
$ws.Range("H67").Value = 28860348
$ws.Range("J67").Value = 62504380
$ws.Range("L67").Value = 62504380
$ws.Range("N67").Value = -62506096
$ws.Range("H74").Value = 4725
$ws.Range("I74").Value = 4400
$ws.Range("K74").Value = 4400
$ws.Range("M74").Value = -3464
$ws.Range("H77").Value = 4725
$ws.Range("I77").Value = 4400
$ws.Range("K77").Value = 22000
$ws.Range("M77").Value = -17320
$ws.Range("H107").Value = 308.17648
$ws.Range("I107").Value = 302.85715
$ws.Range("K107").Value = 302.85715
$ws.Range("M107").Value = 1617.14285
$ws.Range("H137").Value = 15610844
$ws.Range("I137").Value = 1667665.9
$ws.Range("K137").Value = 5002997.699999999
$ws.Range("M137").Value = -5000447.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13938.493
$ws.Range("I32").Value = 14981.981
$ws.Range("J32").Value = 10181.934
$ws.Range("K32").Value = 14981.981
$ws.Range("L32").Value = 10181.934
$ws.Range("M32").Value = -14694.981
$ws.Range("N32").Value = -10755.934
$ws.Range("H45").Value = 3762
$ws.Range("I45").Value = 3585.5293
$ws.Range("K45").Value = 3585.5293
$ws.Range("M45").Value = -3208.5293
$ws.Range("H61").Value = 17308
$ws.Range("I61").Value = 24671.143
$ws.Range("K61").Value = 24671.143
$ws.Range("M61").Value = -24459.143
$ws.Range("H74").Value = 1224.5834
$ws.Range("I74").Value = 923.75
$ws.Range("J74").Value = 1375
$ws.Range("K74").Value = 923.75
$ws.Range("L74").Value = 1375
$ws.Range("M74").Value = -49.75
$ws.Range("N74").Value = -3123
$ws.Range("H77").Value = 1224.5834
$ws.Range("I77").Value = 923.75
$ws.Range("J77").Value = 1375
$ws.Range("K77").Value = 4618.75
$ws.Range("L77").Value = 6875
$ws.Range("M77").Value = -250.75
$ws.Range("N77").Value = -15611
$ws.Range("H122").Value = 3245.111
$ws.Range("I122").Value = 2554.8918
$ws.Range("J122").Value = 6437.375
$ws.Range("K122").Value = 7664.6754
$ws.Range("L122").Value = 19312.125
$ws.Range("M122").Value = -5214.6754
$ws.Range("N122").Value = -24212.125
$ws.Range("H136").Value = 17308
$ws.Range("I136").Value = 24671.143
$ws.Range("K136").Value = 74013.429
$ws.Range("M136").Value = -71463.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 569.4286
$ws.Range("I4").Value = 569.4286
$ws.Range("K4").Value = 569.4286
$ws.Range("M4").Value = -454.4286
$ws.Range("H20").Value = 1819.8966
$ws.Range("J20").Value = 1139.3636
$ws.Range("L20").Value = 1139.3636
$ws.Range("N20").Value = -1633.3636
$ws.Range("H86").Value = 1655
$ws.Range("I86").Value = 1770.7142
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 1770.7142
$ws.Range("L86").Value = 1250
$ws.Range("M86").Value = -647.7141999999999
$ws.Range("N86").Value = -3496
$ws.Range("H89").Value = 1655
$ws.Range("I89").Value = 1770.7142
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 8853.571
$ws.Range("L89").Value = 6250
$ws.Range("M89").Value = -3237.571
$ws.Range("N89").Value = -17482
$ws.Range("H99").Value = 1737644.5
$ws.Range("I99").Value = 2316237
$ws.Range("K99").Value = 2316237
$ws.Range("M99").Value = -2314739

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 4000
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 7000
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 21000
$ws.Range("M21").Value = -2827
$ws.Range("N21").Value = -21346
$ws.Range("H38").Value = 273.3889
$ws.Range("I38").Value = 270.16666
$ws.Range("J38").Value = 275
$ws.Range("K38").Value = 810.4999799999999
$ws.Range("L38").Value = 825
$ws.Range("M38").Value = -463.4999799999999
$ws.Range("N38").Value = -1519
$ws.Range("H56").Value = 47625444
$ws.Range("I56").Value = 47625444
$ws.Range("K56").Value = 47625444
$ws.Range("M56").Value = -47624914
$ws.Range("H98").Value = 1090.1765
$ws.Range("J98").Value = 1141.125
$ws.Range("L98").Value = 3423.375
$ws.Range("N98").Value = -6419.375
$ws.Range("H109").Value = 8981.75
$ws.Range("I109").Value = 3463.5
$ws.Range("K109").Value = 10390.5
$ws.Range("M109").Value = -9350.5
$ws.Range("H121").Value = 1220.6666
$ws.Range("J121").Value = 1220.6666
$ws.Range("L121").Value = 3661.9998
$ws.Range("N121").Value = -6281.9998
$ws.Range("H122").Value = 1073.8
$ws.Range("I122").Value = 793
$ws.Range("J122").Value = 1495
$ws.Range("K122").Value = 7137
$ws.Range("L122").Value = 13455
$ws.Range("M122").Value = -4687
$ws.Range("N122").Value = -18355
$ws.Range("H137").Value = 77284790
$ws.Range("I137").Value = 107144380
$ws.Range("K137").Value = 321433140
$ws.Range("M137").Value = -321428040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H122").Value = 346811.44
$ws.Range("I122").Value = 502476.22
$ws.Range("J122").Value = 4348.9
$ws.Range("K122").Value = 1507428.66
$ws.Range("L122").Value = 13046.7
$ws.Range("M122").Value = -1504978.66
$ws.Range("N122").Value = -17946.7
$ws.Range("H132").Value = 62345.457
$ws.Range("I132").Value = 95482.63
$ws.Range("J132").Value = 6267.154
$ws.Range("K132").Value = 286447.89
$ws.Range("L132").Value = 18801.462
$ws.Range("M132").Value = -283917.89
$ws.Range("N132").Value = -23861.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 883
$ws.Range("I22").Value = 574.5
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 574.5
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -279.5
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 883
$ws.Range("I27").Value = 574.5
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 574.5
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -467.5
$ws.Range("N27").Value = -1714
$ws.Range("H100").Value = 3456.1365
$ws.Range("I100").Value = 3528
$ws.Range("J100").Value = 3369.9
$ws.Range("K100").Value = 3528
$ws.Range("L100").Value = 3369.9
$ws.Range("M100").Value = -2987
$ws.Range("N100").Value = -4451.9
$ws.Range("H132").Value = 3224.484
$ws.Range("I132").Value = 2515.7659
$ws.Range("K132").Value = 7547.297699999999
$ws.Range("M132").Value = -5017.297699999999
$ws.Range("H136").Value = 4235.9287
$ws.Range("I136").Value = 3027.5454
$ws.Range("K136").Value = 9082.636200000001
$ws.Range("M136").Value = -6532.636200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 24975
$ws.Range("J37").Value = 24975
$ws.Range("L37").Value = 24975
$ws.Range("N37").Value = -25381
$ws.Range("H46").Value = 87929
$ws.Range("J46").Value = 87929
$ws.Range("L46").Value = 87929
$ws.Range("N46").Value = -88391
$ws.Range("H122").Value = 4124.608
$ws.Range("I122").Value = 3965.625
$ws.Range("J122").Value = 4392.3687
$ws.Range("K122").Value = 11896.875
$ws.Range("L122").Value = 13177.1061
$ws.Range("M122").Value = -9446.875
$ws.Range("N122").Value = -18077.1061
$ws.Range("H126").Value = 2900.1
$ws.Range("I126").Value = 2382.8572
$ws.Range("J126").Value = 4107
$ws.Range("K126").Value = 7148.571599999999
$ws.Range("L126").Value = 12321
$ws.Range("M126").Value = -4678.571599999999
$ws.Range("N126").Value = -17261
$ws.Range("H134").Value = 87929
$ws.Range("J134").Value = 87929
$ws.Range("L134").Value = 263787
$ws.Range("N134").Value = -268857
$ws.Range("H136").Value = 8245.642
$ws.Range("I136").Value = 2106.375
$ws.Range("J136").Value = 11672.209
$ws.Range("K136").Value = 6319.125
$ws.Range("L136").Value = 35016.627
$ws.Range("M136").Value = -3769.125
$ws.Range("N136").Value = -40116.627
